$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47 ---------------------------------------------------------------
# D47 gets a new, multi-line note (must be written before D46's new string so
# the shared-string table picks up the same insertion order as the authored
# workbook: new index 75 = D47's text, new index 76 = D46's text).
$ws.Range("D47").Value = "課程 : 元辰燈科儀演練`r`n臉書直播...台南道場幹部群組"
$ws.Range("D47").WrapText = $true
$ws.Rows(47).RowHeight = 32.75

# B47 is assigned to 妙一師姐 and loses its "unfilled" yellow highlight,
# picking up the plain bordered style that C47 already uses.
$ws.Range("B47").Value = "妙一師姐"
$ws.Range("C47").Copy()
$ws.Range("B47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 46 -----------------------------------------------------------------
$ws.Range("D46").Value = "仙佛開示"

# --- Row 49 -------------------------------------------------------------
$ws.Range("B49").Value = "玄聖上人師尊"

# --- View state -----------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D47").Select()
